$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Rewrite the Introduction paragraph (paragraph 3) with new text,
#    then insert two new paragraphs after it (inheriting the same
#    ListParagraph / 1440-twip-indent formatting), and finally apply
#    that same paragraph formatting to what was the following blank
#    paragraph.
# ------------------------------------------------------------------

$introPara = $d.Paragraphs.Item(3)
$introRange = $introPara.Range
$introRange.MoveEnd(1, -1) | Out-Null
$introRange.Text = "Malaria, a deadly vector-borne disease caused by Plasmodium parasites, remains a significant public health concern in many regions across the world. Among these regions, Africa bears the heaviest burden, with a disproportionate share of the global malaria cases and deaths."

# Insert "Efforts to combat malaria..." right after paragraph 3.
$introPara = $d.Paragraphs.Item(3)
$introPara.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range
$r4.MoveEnd(1, -1) | Out-Null
$r4.Text = "Efforts to combat malaria in Africa encompass a range of strategies, from vector control (using insecticide treated bed nets) to preventative treatments for high-risk groups (Intermittent Preventative treatment for pregnant women). "

# Insert "In this project, we delve..." right after paragraph 4.
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item(5)
$r5 = $p5.Range
$r5.MoveEnd(1, -1) | Out-Null
$r5.Text = "In this project, we delve into the multifaceted issue of malaria in Africa, exploring its epidemiology, prevention strategies, and ongoing efforts to combat this pervasive threat."

# The paragraph that used to immediately follow the introduction (a
# blank paragraph with no special style) is now paragraph 6. Give it
# the same ListParagraph style / 1440-twip left indent as its
# neighbours (preserving its paragraph-mark language formatting).
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Style = "List Paragraph"
$p6.Range.ParagraphFormat.LeftIndent = 72
$p6.Range.LanguageID = "en-US"

# ------------------------------------------------------------------
# 2. Insert a new bold "Methodology" heading right before "Results".
# ------------------------------------------------------------------

$resultsIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Results") {
        $resultsIndex = $i
        break
    }
}

$resultsPara = $d.Paragraphs.Item($resultsIndex)
$resultsPara.Range.InsertParagraphBefore() | Out-Null

# The freshly inserted (blank) paragraph now occupies $resultsIndex,
# having inherited "Results"' bold / numbered-heading formatting; we
# just need to give it its text.
$methodPara = $d.Paragraphs.Item($resultsIndex)
$mr = $methodPara.Range
$mr.MoveEnd(1, -1) | Out-Null
$mr.Text = "Methodology"

# ------------------------------------------------------------------
# 3. "Conclusion" -> "Conclusion."
# ------------------------------------------------------------------

$d.Content.Find.Execute("Conclusion", $true, $true, $false, $false, $false, $true, 1, $false, "Conclusion.", 2) | Out-Null
